# Apply updated "想去人数" (interest count) figures and one ticket-status
# change, mirroring the gh-pages data refresh described in the commit.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 290
$ws1.Range("F3").Value  = 661
$ws1.Range("F10").Value = 593
$ws1.Range("F14").Value = 6046
$ws1.Range("F15").Value = 651
$ws1.Range("F17").Value = 20
$ws1.Range("F18").Value = 251
$ws1.Range("F19").Value = 176
$ws1.Range("F21").Value = 578
$ws1.Range("F22").Value = 14
$ws1.Range("F23").Value = 50
$ws1.Range("F25").Value = 153
$ws1.Range("F26").Value = 1335
$ws1.Range("F28").Value = 1019
$ws1.Range("F29").Value = 62
$ws1.Range("F30").Value = 2081
$ws1.Range("F31").Value = 194
$ws1.Range("F32").Value = 24
$ws1.Range("F35").Value = 3347

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("G3").Value  = "已停售"
$ws2.Range("F7").Value  = 99
$ws2.Range("F11").Value = 652
$ws2.Range("F16").Value = 79
$ws2.Range("F21").Value = 358
$ws2.Range("F23").Value = 4060
$ws2.Range("F27").Value = 159
$ws2.Range("F28").Value = 220
$ws2.Range("F29").Value = 73

# --- Sheet "本地生活" (local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value  = 2605
$ws3.Range("F6").Value  = 1153
$ws3.Range("F8").Value  = 1505
$ws3.Range("F9").Value  = 421
$ws3.Range("F12").Value = 672

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 2605
$ws4.Range("F5").Value  = 1153
$ws4.Range("F6").Value  = 1505
$ws4.Range("F7").Value  = 421
$ws4.Range("F9").Value  = 290
$ws4.Range("F10").Value = 661
$ws4.Range("F13").Value = 672
$ws4.Range("F14").Value = 593
$ws4.Range("F15").Value = 99
$ws4.Range("F19").Value = 6046
$ws4.Range("F21").Value = 651
$ws4.Range("F23").Value = 251
$ws4.Range("F24").Value = 176
$ws4.Range("F26").Value = 578
$ws4.Range("F28").Value = 79
$ws4.Range("F31").Value = 14
$ws4.Range("F34").Value = 358
$ws4.Range("F38").Value = 159
$ws4.Range("F39").Value = 220
$ws4.Range("F40").Value = 62
$ws4.Range("F41").Value = 73
$ws4.Range("F43").Value = 2081
$ws4.Range("F46").Value = 194
$ws4.Range("F47").Value = 24
$ws4.Range("F49").Value = 3347
